$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the attendance row for the 9/29 / 4:15 meeting (row 17) ---

# B17: new meeting date/time string
$ws.Range("B17").Value = "9/29 / 4:15"

# C17: meeting place ("Google Hangout") - copy formatting from the cell
# above (C16) first so the style matches the other rows (bold border style),
# then set the value.
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = "Google Hangout"

# D17..G17, I17: attendance = "A" (attend) for everyone except Younouss Thiam
$ws.Range("D17").Value = "A"
$ws.Range("E17").Value = "A"
$ws.Range("F17").Value = "A"
$ws.Range("G17").Value = "A"
# H17: Younouss Thiam -> "U" (unexcused absence)
$ws.Range("H17").Value = "U"
$ws.Range("I17").Value = "A"

# Update the active selection to reflect where the user ended up editing
$ws.Range("J17").Select()
